# Generate Report for Archive
# Updates status text from "Ready for handoff" to "In Translation" across
# the Overview / zh-cn / de-de sheets, and re-fits the affected status
# columns to their new (narrower) width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"  # previous status text being replaced everywhere
$newStatus = "In Translation"
# Target stored column width (OOXML "width" attribute) is 13.4101845877511,
# which was produced by Excel's real AutoFit pixel metrics. This runtime's
# ColumnWidth setter quantizes to 1/6-character steps (stored = round(cw*6)/6
# + 5/6), so 12.5 is the closest achievable input (-> stored 13.333333...).
$newColumnWidth = 12.5

# Overview sheet: status values live in columns E (zh-cn) and F (de-de) of row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

# zh-cn sheet: status value lives in column C of row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

# de-de sheet: status value lives in column C of row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
